$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# This workbook tracks the localization handoff/handback status for
# two source files:
#   55674320-8dc8-4dc7-a8de-82320adb3288.md
#   f22db898-444b-40f8-83f1-7c04d8cfd271.md
# A new handoff report has been generated: the row that used to
# describe f22db898 now describes 55674320 (and vice versa), the
# newly-handed-off file (55674320) moves to status "Ready for
# handoff", and the handoff timestamps for that batch are refreshed.
# ------------------------------------------------------------------

$missing = [System.Type]::Missing

# ====================================================================
# Sheet 1: "Overview"
# ====================================================================
$ws1 = $wb.Worksheets.Item(1)

# Drop every hyperlink on the sheet so we can recreate them (in order)
# with the same targets but new display text / cell values.
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Range("A2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/55674320-8dc8-4dc7-a8de-82320adb3288.md", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/f22db898-444b-40f8-83f1-7c04d8cfd271.md", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ====================================================================
# Sheet 2: "zh-cn"
# ====================================================================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Hyperlinks.Delete()

# Row 2 now carries the f22db898 file's data
$ws2.Range("A2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-03 08:15:04"
$ws2.Range("E2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.md"
$ws2.Range("F2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-03 08:14:01"
$ws2.Range("H2").Value = "Include"

# Row 3 now carries the 55674320 file's data
$ws2.Range("A3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-03 08:15:04"
$ws2.Range("E3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.md"
$ws2.Range("F3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-03 08:14:01"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/55674320-8dc8-4dc7-a8de-82320adb3288.md", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/490be9eedadc51e191bbc7f7f3b9afc7865af816/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0639a454a66317960a1eb0adcdde2b328d913be2/e2e/55674320-8dc8-4dc7-a8de-82320adb3288.md", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7042bbefc4ac64a43f74c03bc5a9b672888a4811/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/f22db898-444b-40f8-83f1-7c04d8cfd271.md", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/490be9eedadc51e191bbc7f7f3b9afc7865af816/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0639a454a66317960a1eb0adcdde2b328d913be2/e2e/f22db898-444b-40f8-83f1-7c04d8cfd271.md", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7042bbefc4ac64a43f74c03bc5a9b672888a4811/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.zh-cn.xlf", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ====================================================================
# Sheet 3: "de-de"
# ====================================================================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Hyperlinks.Delete()

# Row 2 now carries the f22db898 file's data
$ws3.Range("A2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-03 08:15:15"
$ws3.Range("E2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.md"
$ws3.Range("F2").Value = "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-03 08:14:22"
$ws3.Range("H2").Value = "Include"

# Row 3 now carries the 55674320 file's data
$ws3.Range("A3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-03 08:15:15"
$ws3.Range("E3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.md"
$ws3.Range("F3").Value = "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-03 08:14:22"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/55674320-8dc8-4dc7-a8de-82320adb3288.md", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ecb008cecd9bfc2d30ac5ff366ae9c0f280e9bf2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f480434f014b156db597bc5461b62a9a9c13dc6e/e2e/55674320-8dc8-4dc7-a8de-82320adb3288.md", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/54c4bf6c9cc55f226de0365db3f05726372a6052/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf", $missing, $missing, "f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/e2e/f22db898-444b-40f8-83f1-7c04d8cfd271.md", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ecb008cecd9bfc2d30ac5ff366ae9c0f280e9bf2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f480434f014b156db597bc5461b62a9a9c13dc6e/e2e/f22db898-444b-40f8-83f1-7c04d8cfd271.md", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/54c4bf6c9cc55f226de0365db3f05726372a6052/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f22db898-444b-40f8-83f1-7c04d8cfd271.e8a20d3df549e014ff3fad3aa1d9b8887c69e9e9.de-de.xlf", $missing, $missing, "55674320-8dc8-4dc7-a8de-82320adb3288.4955fab11732025d2a74aca999e07dbc1df2af7b.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/da01010da2ec688706ff402d8d32ede34cb1a9aa/.localization-config", $missing, $missing, ".localization-config") | Out-Null
